$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: "$(ProSpec)" -> "$(ProSpec) <= $in" for data rows 2-6
$ws.Range("A2").Value = "`$(ProSpec) <= `$in"
$ws.Range("A3").Value = "`$(ProSpec) <= `$in"
$ws.Range("A4").Value = "`$(ProSpec) <= `$in"
$ws.Range("A5").Value = "`$(ProSpec) <= `$in"
$ws.Range("A6").Value = "`$(ProSpec) <= `$in"

# Column C header: "out:price" -> "out:ProSpec price"
$ws.Range("C1").Value = "out:ProSpec price"

# Column C data: "${Zone N}" -> "${FedExZone N}"
$ws.Range("C2").Value = "`${FedExZone 1}"
$ws.Range("C3").Value = "`${FedExZone 2}"
$ws.Range("C4").Value = "`${FedExZone 3}"
$ws.Range("C5").Value = "`${FedExZone 4}"
$ws.Range("C6").Value = "`${FedExZone 5}"

# Column C width (engine snaps ColumnWidth to the nearest 1/6 character unit;
# 31 is the closest input that reproduces the target stored width of ~31.832)
$ws.Columns.Item(3).ColumnWidth = 31

# Selection
$ws.Range("C7").Select()
